$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D..AJ across rows 2..9 (row => column-letter => value)
$data = @{
    2 = @{ D=195723; E=10031; F=10031; G=7973; H=5086; I=2092; J=2995; K=230393; L=139694; M=90699; N=34054; O=56645; P=1575; Q=17610; R=-19327; S=-735; T=10108; U=7502; V=89971; W=5.13; X=2.6; Y=6.31; Z=2.25; AA=154.02; AB=2124.62; AC=5884; AD=25; AE=104465; AF=1.41; AG=950; AH=0.65; AJ=29084427 }
    3 = @{ D=211667; E=12253; F=12253; G=7349; H=5511; I=2058; J=3453; K=235621; L=136987; M=98635; N=36012; O=62623; P=1577; Q=22956; R=-14615; S=-6371; T=11280; U=11676; V=86160; W=5.79; X=2.6; Y=5.87; Z=2.37; AA=138.88; AB=2228.07; AC=5781; AD=40.94; AE=110305; AF=2.15; AG=1350; AH=0.57; AJ=29133748 }
    4 = @{ D=239542; E=12529; F=12529; G=7953; H=5698; I=2230; J=3468; K=270096; L=157551; M=112544; N=37358; O=75186; P=1579; Q=20917; R=-33447; S=12305; T=13744; U=7174; V=97963; W=5.23; X=2.38; Y=6.08; Z=2.25; AA=139.99; AB=2338.27; AC=6257; AD=28.18; AE=114276; AF=1.54; AG=1350; AH=0.77; AJ=29176998 }
    5 = @{ D=268986; E=13260; F=13260; G=15752; H=11377; I=4637; J=6740; K=293254; L=173332; M=119922; N=39298; O=80624; P=1579; Q=22404; R=-28638; S=7182; T=22729; U=-325; V=105978; W=4.93; X=4.23; Y=12.1; Z=4.04; AA=144.54; AB=2601.61; AC=13002; AD=13.16; AE=120209; AF=1.42; AG=1450; AH=0.85; AJ=29176998 }
    6 = @{ D=295234; E=13325; F=13325; G=13197; H=8800; I=2761; K=317295; L=192077; M=125218; N=40587; P=1579; Q=15411; R=-21260; S=8274; T=24583; U=-9172; V=120449; W=4.51; X=2.98; Y=6.91; Z=2.88; AA=153.39; AB=2741.67; AC=7741; AD=15.7; AE=125257; AF=0.97; AJ=29176998 }
    7 = @{ D=338597; E=14584; G=8494; H=5914; I=3210; K=385921; L=248418; M=137503; N=43851; P=1790; Q=36522; R=-45096; S=13524; T=22006; U=-695; W=4.31; X=1.75; Y=7.6; Z=1.68; AA=180.66; AC=9001; AD=9.2; AE=135384; AF=0.61; AG=1481; AH=1.79; AI=13.46 }
    8 = @{ D=362502; E=16567; G=9925; H=6972; I=3250; K=397002; L=252433; M=144571; N=46725; P=1790; Q=31650; R=-16891; S=-4332; T=14215; U=14952; W=4.57; X=1.92; Y=7.17; Z=1.78; AA=174.61; AC=9112; AD=9.09; AE=144256; AF=0.57; AG=1506; AH=1.82; AI=13.52 }
    9 = @{ D=385702; E=17777; G=11374; H=7998; I=3715; K=408729; L=256093; M=152636; N=50029; P=1790; Q=32391; R=-16332; S=-4595; T=12893; U=16206; W=4.61; X=2.07; Y=7.68; Z=1.98; AA=167.78; AC=10418; AD=7.95; AE=154457; AF=0.54; AG=1531; AH=1.85; AI=12.03 }
}

# Columns that must become empty (were present before, removed in target)
$clearCells = @("AI2", "AI3", "AI4", "AI5", "AG6", "AH6", "AI6")

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}
